$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ValidUserLogin")

# Update C1 on the first sheet to reference the new "runMode" string
# (was "RunMode", now lowercase "runMode")
$ws.Range("C1").Value = "runMode"

# Update the active selection on the first sheet to D2
$ws.Range("D2").Select()
